# Insert a new data row at row 419 (a new weekly price observation),
# shifting all existing rows from 419..545 down to 420..546.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("419:419").Insert()

# Populate the newly inserted row 419 with the new observation.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Categoria ID,
# G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo, L Precio maximo,
# M Precio promedio ponderado, N Unidad de comercializacion, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificacion
$ws.Cells.Item(419, 1).Value = 4
$ws.Cells.Item(419, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(419, 3).Value = "Los Lagos"
$ws.Cells.Item(419, 4).Value = "2023-10-16"
$ws.Cells.Item(419, 5).Value = 10
$ws.Cells.Item(419, 6).Value = 100114014
$ws.Cells.Item(419, 7).Value = "Betarraga"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 250
$ws.Cells.Item(419, 11).Value = 1000
$ws.Cells.Item(419, 12).Value = 1000
$ws.Cells.Item(419, 13).Value = 1000
$ws.Cells.Item(419, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(419, 15).Value = "Región Metropolitana"
$ws.Cells.Item(419, 16).Value = 200
$ws.Cells.Item(419, 17).Value = 5
$ws.Cells.Item(419, 18).Value = "Hortaliza"
